# Apply the change: insert a new inventory-shortage row for
# "ATOR 10MG 7 TAB." as item #4 (between ANTINAL and BETADERM),
# shifting every following item down by one position, add a brand
# new trailing row for item #15 ("سرنجات 5 سم"), and bump the
# totals row accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Make room for the new trailing item row by inserting a blank
#    row at 21 (this pushes the old totals row 21->22 and the old
#    footer row 22->23, exactly like the target workbook).
# ---------------------------------------------------------------
$ws.Rows.Item(21).Insert()

# Give the new row the same look as the other item rows (copy the
# formatting from row 20, which is the row directly above it).
$fmtSrc = $ws.Range("A20:Q20")
$fmtDst = $ws.Range("A21:Q21")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = 0

$ws.Rows.Item(21).RowHeight = 25.5

$ws.Range("A21:B21").Merge()
$ws.Range("C21:G21").Merge()
$ws.Range("H21:K21").Merge()
$ws.Range("L21:M21").Merge()
$ws.Range("N21:O21").Merge()

# ---------------------------------------------------------------
# 2. Rewrite the item list, rows 10-21, shifting every item from
#    row N (old) down into row N+1 (new) and writing the brand new
#    "ATOR 10MG 7 TAB." item into row 10.
#
#    Walk rows 20 down to 10, copying each row's data into the row
#    below it (row+1), so the list for rows 11-21 becomes the old
#    rows 10-20.
# ---------------------------------------------------------------
for ($r = 20; $r -ge 10; $r--) {
    $num        = $ws.Cells.Item($r, 1).Value2
    $name       = $ws.Cells.Item($r, 3).Value2
    $balance    = $ws.Cells.Item($r, 8).Value2
    $orderLimit = $ws.Cells.Item($r, 12).Value2
    $price      = $ws.Cells.Item($r, 14).Value2
    $salePrice  = $ws.Cells.Item($r, 16).Value2
    $trans      = $ws.Cells.Item($r, 17).Value2

    $newNum = $num + 1
    $dstRow = $r + 1

    $ws.Cells.Item($dstRow, 1).Value = $newNum
    $ws.Cells.Item($dstRow, 3).Value = $name
    $ws.Cells.Item($dstRow, 8).Value = $balance
    $ws.Cells.Item($dstRow, 12).Value = $orderLimit
    $ws.Cells.Item($dstRow, 14).Value = $price
    $ws.Cells.Item($dstRow, 16).Value = $salePrice
    $ws.Cells.Item($dstRow, 17).Value = $trans
}

# Now write the new item into row 10.
$ws.Cells.Item(10, 1).Value = 4
$ws.Cells.Item(10, 3).Value = "ATOR 10MG 7 TAB."
$ws.Cells.Item(10, 8).Value = "0:0"
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 14).Value = "45.00"
$ws.Cells.Item(10, 16).Value = "90.0000"
$ws.Cells.Item(10, 17).Value = "2:0"

# ---------------------------------------------------------------
# 3. Update the totals row (now row 22): add the new item's sale
#    price total (90.0000) to the previous grand total.
# ---------------------------------------------------------------
$total = $ws.Cells.Item(22, 14).Value2
$newTotal = $total + 90
$ws.Cells.Item(22, 14).Value = $newTotal
